$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "(0.0)"
$ws.Range("C4").Value = "(0.04)"
$ws.Range("D4").Value = "(0.15)"
$ws.Range("E4").Value = "(0.03)"
$ws.Range("F4").Value = "(0.03)"
$ws.Range("G4").Value = "(0.06)"

$ws.Range("B6").Value = "(0.0)"
$ws.Range("C6").Value = "(0.09)"
$ws.Range("D6").Value = "(0.06)"
$ws.Range("E6").Value = "(0.06)"
$ws.Range("F6").Value = "(0.03)"
$ws.Range("G6").Value = "(0.06)"
